$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster = ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.5587383333333333
$ws.Range("H2").Value = 1.676215
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 79.66420866666668
$ws.Range("N2").Value = 238.992626
$ws.Range("O2").Value = 0.1048332405251988
$ws.Range("P2").Value = 0.1104048582881303
$ws.Range("Q2").Value = 44.51144717673223
$ws.Range("R2").Value = 400.6030245905901
$ws.Range("S2").Value = 0.1048332405251988
$ws.Range("T2").Value = 0.1104048582881303

# Row 3 (Target cluster = FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.5587383333333333
$ws.Range("H3").Value = 1.676215
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 316.1112466666667
$ws.Range("N3").Value = 948.33374
$ws.Range("O3").Value = 0.4159831235277584
$ws.Range("P3").Value = 0.438091559254019
$ws.Range("Q3").Value = 176.6234711104555
$ws.Range("R3").Value = 1589.6112399941
$ws.Range("S3").Value = 0.4159831235277584
$ws.Range("T3").Value = 0.438091559254019

# Row 4 (Target cluster = M1)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.5587383333333333
$ws.Range("H4").Value = 1.676215
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 116.842289
$ws.Range("N4").Value = 350.526867
$ws.Range("O4").Value = 0.1537573270514019
$ws.Range("P4").Value = 0.1619291344885147
$ws.Range("Q4").Value = 65.28426581871167
$ws.Range("R4").Value = 587.5583923684051
$ws.Range("S4").Value = 0.1537573270514019
$ws.Range("T4").Value = 0.1619291344885147

# Row 5 (Target cluster = M2)
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.5587383333333333
$ws.Range("H5").Value = 1.676215
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 132.2480086666667
$ws.Range("N5").Value = 396.744026
$ws.Range("O5").Value = 0.1740303146616488
$ws.Range("P5").Value = 0.1832795793757766
$ws.Range("Q5").Value = 73.89203194906555
$ws.Range("R5").Value = 665.0282875415901
$ws.Range("S5").Value = 0.1740303146616488
$ws.Range("T5").Value = 0.1832795793757766

# Row 6 (Target cluster = sCs)
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.5587383333333333
$ws.Range("H6").Value = 1.676215
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 115.04788
$ws.Range("N6").Value = 230.09576
$ws.Range("O6").Value = 0.1513959942339921
$ws.Range("P6").Value = 0.1062948685935592
$ws.Range("Q6").Value = 64.28166072473333
$ws.Range("R6").Value = 385.6899643484
$ws.Range("S6").Value = 0.1513959942339921
$ws.Range("T6").Value = 0.1062948685935592
